$d = $word.ActiveDocument
$d.Content.Find.Execute("Line(191 to 195)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Line(268 to 271)", 2)
